$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new blank row at position 35, pushing the "farm_address_comments"
# row (and everything after it) down by one.
$ws.Rows.Item(35).Insert()

# Row 34 used to be the single "farm_location" geopoint field (hint +
# appearance "placement-map"). It now becomes the latitude half of the
# split field: drop the hint/appearance, add a numeric range constraint,
# constraint message and calculation that pulls the latitude.
$ws.Range("D34").ClearContents()
$ws.Range("L34").ClearContents()

$ws.Range("C34").Value = "Latitud de la finca"
$ws.Range("G34").Value = '(. >= -90) and (. <= 90)'
$ws.Range("H34").Value = "La latitud debe ser un campo numérico entre -90 y 90"
$ws.Range("J34").Value = 'number(pulldata(''farm'',''lat'', ''farm_key'', ${farm_id}))'

# New row 35 is the longitude half of the split field.
$ws.Range("A35").Value = "decimal"
$ws.Range("B35").Value = "farm_location_lon"
$ws.Range("C35").Value = "Longitud de la finca"
$ws.Range("E35").Value = "yes"
$ws.Range("G35").Value = '(. >= -180) and (. <= 180)'
$ws.Range("H35").Value = "La longitud debe ser un campo numérico entre -180 y 180"
$ws.Range("J35").Value = 'number(pulldata(''farm'',''lon'', ''farm_key'', ${farm_id}))'

$ws.Range("C35").Select()

# The survey sheet used to span A1:L195; after inserting the new row it
# spans A1:L196, so the (hidden) filter-database defined name needs to grow
# along with it.
$n = $wb.Names.Item("_xlnm._FilterDatabase")
$n.RefersTo = "=survey!`$A`$1:`$L`$196"
